$d = $word.ActiveDocument

$replacements = @(
    @("2024-01-29 Monday", "2024-01-30 Tuesday"),
    @("66×54=3564", "26×33=858"),
    @("70×52=3640", "20×94=1880"),
    @("48×68=3264", "71×20=1420"),
    @("60×57=3420", "86×61=5246"),
    @("88×94=8272", "63×21=1323"),
    @("67×71=4757", "37×72=2664"),
    @("75×40=3000", "73×56=4088"),
    @("67×89=5963", "95×60=5700"),
    @("31×64=1984", "22×70=1540"),
    @("53×31=1643", "72×30=2160"),
    @("80×62=4960", "99×12=1188"),
    @("37×32=1184", "21×40=840"),
    @("50×18=900", "93×70=6510"),
    @("13×30=390", "96×92=8832"),
    @("30×17=510", "42×25=1050"),
    @("54×19=1026", "43×79=3397"),
    @("74×24=1776", "32×40=1280"),
    @("44×69=3036", "30×46=1380"),
    @("56×39=2184", "48×72=3456"),
    @("24×24=576", "15×69=1035"),
    @("19×24=456", "68×48=3264"),
    @("20×99=1980", "49×29=1421"),
    @("20×74=1480", "33×37=1221"),
    @("32×22=704", "62×35=2170"),
    @("30×71=2130", "40×79=3160")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
